# Column K ("Fecha") currently holds the date "07 08 24" as inline-string
# text in every data row. Replace it with the real date 2024-08-07 stored
# as a numeric (serial) value, formatted as a date/time.
#
# Each cell is first given the intermediate format "yyyy-mm-dd h:mm:ss"
# (this is the format a date value naturally picks up) before the real
# date value is written, and only then is the format switched to the
# final display format "YYYY-MM-DD HH:MM:SS". This mirrors how the
# workbook was produced (extracting/parsing the date, then formatting it)
# and keeps both format codes registered in the workbook's number-format
# table, exactly like the source change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValue = Get-Date -Year 2024 -Month 8 -Day 7 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 11)   # column K
    $cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
    $cell.Value = $dateValue
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
